$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value2 = '[Carol H.%Yan%c1yan@health.ucsd.edu%0,                           Farhoud%Faraji%NULL%0,                           Farhoud%Faraji%NULL%0,                           Divya P.%Prajapati%NULL%0,                           Benjamin T.%Ostrander%NULL%0,                           Adam S.%DeConde%NULL%0]'
$ws.Range('E3').Value2 = '[Á.%Beltrán‐Corbellini%beltran_corbellini@hotmail.com%0,                           J. L.%Chico‐García%NULL%0,                           J. L.%Chico‐García%NULL%0,                           J.%Martínez‐Poles%NULL%0,                           F.%Rodríguez‐Jorge%NULL%0,                           E.%Natera‐Villalba%NULL%0,                           J.%Gómez‐Corral%NULL%0,                           J.%Gómez‐Corral%NULL%0,                           A.%Gómez‐López%NULL%0,                           E.%Monreal%NULL%0,                           P.%Parra‐Díaz%NULL%0,                           J. L.%Cortés‐Cuevas%NULL%0,                           J. L.%Cortés‐Cuevas%NULL%0,                           J. C.%Galán%NULL%0,                           C.%Fragola‐Arnau%NULL%0,                           J.%Porta‐Etessam%NULL%0,                           J.%Masjuan%NULL%0,                           A.%Alonso‐Cánovas%NULL%0,                           A.%Alonso‐Cánovas%NULL%0]'
$ws.Range('E4').Value2 = '[François%Bénézit%NULL%0,                           Paul%Le Turnier%NULL%0,                           Charles%Declerck%NULL%0,                           Cécile%Paillé%NULL%0,                           Matthieu%Revest%NULL%0,                           Vincent%Dubée%NULL%0,                           Pierre%Tattevin%pierre.tattevin@chu-rennes.fr%0,                           Cédric%Arvieux%NULL%0,                           Marion%Baldeyrou%NULL%0,                           Jean-Marc%Chapplain%NULL%0,                           Pauline%Comacle%NULL%0,                           Solène%Patrat-Delon%NULL%0,                           Anne%Maillard%NULL%0,                           Mélanie%Poinot%NULL%0,                           Charlotte%Pronier%NULL%0,                           Faouzi%Souala%NULL%0,                           Vincent%Thibault%NULL%0,                           Pierre%Abgueguen%NULL%0,                           Hélène%Cormier%NULL%0,                           Valérie%Delbos%NULL%0,                           Marine%de la Chapelle%NULL%0,                           Alexandra%Ducancelle%NULL%0,                           Rafael%Mahieu%NULL%0,                           Valérie%Rabier%NULL%0,                           Sami%Rehaiem%NULL%0,                           Yves%Vandamme%NULL%0,                           Charlotte%Biron%NULL%0,                           Jeanne%Brochon%NULL%0,                           David%Boutoille%NULL%0,                           Marie%Chauveau%NULL%0,                           Colin%Deschanvres%NULL%0,                           Benjamin J%Gaborit%NULL%0,                           Joël%Jenvrin%NULL%0,                           Raphaël%Lecomte%NULL%0,                           Maeva%Lefebvre%NULL%0,                           François%Raffi%NULL%0]'
$ws.Range('E5').Value2 = '[Andrea%Giacomelli%andrea.giacomelli@unimi.it%0,                           Laura%Pezzati%NULL%0,                           Laura%Pezzati%NULL%0,                           Federico%Conti%NULL%0,                           Dario%Bernacchia%NULL%0,                           Matteo%Siano%NULL%0,                           Letizia%Oreni%NULL%0,                           Stefano%Rusconi%NULL%0,                           Cristina%Gervasoni%NULL%0,                           Anna Lisa%Ridolfo%NULL%0,                           Giuliano%Rizzardini%NULL%0,                           Spinello%Antinori%NULL%0,                           Massimo%Galli%NULL%0,                           Massimo%Galli%NULL%0]'
$ws.Range('E6').Value2 = '[T.%Klopfenstein%NULL%0,                           N.J.%Kadiane-Oussou%NULL%0,                           L.%Toko%NULL%0,                           P.-Y.%Royer%NULL%0,                           Q.%Lepiller%NULL%0,                           V.%Gendrin%NULL%0,                           S.%Zayet%NULL%0]'
$ws.Range('E7').Value2 = '[Jerome R.%Lechien%Jerome.Lechien@umons.ac.be%0,                           Carlos M.%Chiesa-Estomba%NULL%0,                           Daniele R.%De Siati%NULL%0,                           Mihaela%Horoi%NULL%0,                           Serge D.%Le Bon%NULL%0,                           Alexandra%Rodriguez%NULL%0,                           Didier%Dequanter%NULL%0,                           Serge%Blecic%NULL%0,                           Fahd%El Afia%NULL%0,                           Lea%Distinguin%NULL%0,                           Younes%Chekkoury-Idrissi%NULL%0,                           Stéphane%Hans%NULL%0,                           Irene Lopez%Delgado%NULL%0,                           Christian%Calvo-Henriquez%NULL%0,                           Philippe%Lavigne%NULL%0,                           Chiara%Falanga%NULL%0,                           Maria Rosaria%Barillari%NULL%0,                           Giovanni%Cammaroto%NULL%0,                           Mohamad%Khalife%NULL%0,                           Pierre%Leich%NULL%0,                           Christel%Souchay%NULL%0,                           Camelia%Rossi%NULL%0,                           Fabrice%Journe%NULL%0,                           Julien%Hsieh%NULL%0,                           Myriam%Edjlali%NULL%0,                           Robert%Carlier%NULL%0,                           Laurence%Ris%NULL%0,                           Andrea%Lovato%NULL%0,                           Cosimo%De Filippis%NULL%0,                           Frederique%Coppee%NULL%0,                           Nicolas%Fakhry%NULL%0,                           Tareck%Ayad%NULL%0,                           Sven%Saussez%NULL%0]'
$ws.Range('E8').Value2 = '[Ling%Mao%xref no email%0,       Huijuan%Jin%xref no email%0,       Mengdie%Wang%xref no email%0,       Yu%Hu%xref no email%0,       Shengcai%Chen%xref no email%0,       Quanwei%He%xref no email%0,       Jiang%Chang%xref no email%0,       Candong%Hong%xref no email%0,       Yifan%Zhou%xref no email%0,       David%Wang%xref no email%0,       Xiaoping%Miao%xref no email%0,       Yanan%Li%xref no email%0,       Bo%Hu%xref no email%0]'
$ws.Range('E9').Value2 = '[ Cristina%Menni%null%1,                  Ana M.%Valdes%null%1,                  Maxim B.%Freidin%null%1,                  Carole H.%Sudre%null%1,                  Long H.%Nguyen%null%1,                  David A.%Drew%null%1,                  Sajaysurya%Ganesh%null%1,                  Thomas%Varsavsky%null%1,                  M. Jorge%Cardoso%null%1,                  Julia S.%El-Sayed Moustafa%null%1,                  Alessia%Visconti%null%1,                  Pirro%Hysi%null%1,                  Ruth C. E.%Bowyer%null%1,                  Massimo%Mangino%null%1,                  Mario%Falchi%null%1,                  Jonathan%Wolf%null%1,                  Sebastien%Ourselin%null%1,                  Andrew T.%Chan%null%1,                  Claire J.%Steves%null%1,                  Tim D.%Spector%null%1]'
$ws.Range('E10').Value2 = '[Shima T.%Moein%NULL%1,                           Seyed MohammadReza%Hashemian%NULL%1,                           Babak%Mansourafshar%NULL%1,                           Ali%Khorram‐Tousi%NULL%1,                           Payam%Tabarsi%NULL%1,                           Richard L.%Doty%doty@pennmedicine.upenn.edu%1]'
$ws.Range('E11').Value2 = '[Borsetto%Daniele%coreGivesNoEmail%1,                          Boscolo-Rizzo%Paolo%coreGivesNoEmail%1,                          Cazzador%Diego%coreGivesNoEmail%1,                          Fabbris%Cristoforo%coreGivesNoEmail%1,                          Hopkins%Claire%coreGivesNoEmail%1,                          Polesel%Jerry%coreGivesNoEmail%1,                          Spinato%Giacomo%coreGivesNoEmail%1]'
$ws.Range('E12').Value2 = '[Alma%Tostmann%NULL%1,                           John%Bradley%NULL%2,                           John%Bradley%NULL%0,                           Teun%Bousema%NULL%2,                           Teun%Bousema%NULL%0,                           Wing-Kee%Yiek%NULL%2,                           Wing-Kee%Yiek%NULL%0,                           Minke%Holwerda%NULL%1,                           Chantal%Bleeker-Rovers%NULL%1,                           Jaap%ten Oever%NULL%1,                           Corianne%Meijer%NULL%2,                           Corianne%Meijer%NULL%0,                           Janette%Rahamat-Langendoen%NULL%1,                           Joost%Hopman%NULL%2,                           Joost%Hopman%NULL%0,                           Nannet%van der Geest-Blankert%NULL%1,                           Heiman%Wertheim%NULL%1]'
$ws.Range('E13').Value2 = '[Luigi Angelo%Vaira%luigi.vaira@gmail.com%0,                          Giovanna%Deiana%NULL%0,                          Giovanna%Deiana%NULL%0,                          Alessandro Giuseppe%Fois%NULL%1,                          Pietro%Pirina%NULL%0,                          Giordano%Madeddu%NULL%0,                          Andrea%De Vito%NULL%0,                          Sergio%Babudieri%NULL%0,                          Marzia%Petrocelli%NULL%0,                          Antonello%Serra%NULL%0,                          Francesco%Bussu%NULL%0,                          Enrica%Ligas%NULL%1,                          Giovanni%Salzano%NULL%0,                          Giacomo%De Riu%NULL%0]'
$ws.Range('E14').Value2 = '[Liang En%Wee%ian.wee@mohh.com.sg%1,                          Yvonne Fu Zi%Chan%NULL%2,                          Yvonne Fu Zi%Chan%NULL%0,                          Neville Wei Yang%Teo%NULL%1,                          Benjamin Pei Zhi%Cherng%NULL%1,                          Siew Yee%Thien%NULL%1,                          Hei Man%Wong%NULL%1,                          Limin%Wijaya%NULL%1,                          Song Tar%Toh%NULL%1,                          Thuan Tong%Tan%NULL%1]'
$ws.Range('E15').Value2 = '[Carol H.%Yan%c1yan@health.ucsd.edu%0,                           Farhoud%Faraji%NULL%0,                           Farhoud%Faraji%NULL%0,                           Divya P.%Prajapati%NULL%0,                           Divya P.%Prajapati%NULL%0,                           Christine E.%Boone%NULL%2,                           Christine E.%Boone%NULL%0,                           Adam S.%DeConde%NULL%0,                           Adam S.%DeConde%NULL%0]'
